$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 24.60000000000041
$ws.Range("H2").Value = [double]"6.05073694323055e-05"
$ws.Range("I2").Value = [double]"6.05073694323055e-05"
$ws.Range("L2").Value = 58.07070347889538
$ws.Range("M2").Value = "[31.31375803263937, 84.8276489251514]"
$ws.Range("N2").Value = [double]"7.210756288578146e-05"
$ws.Range("O2").Value = [double]"7.210756288578146e-05"
$ws.Range("P2").Value = 1.566079220708426
$ws.Range("Q2").Value = "[0.9748685912040411, 2.157289850212811]"
$ws.Range("R2").Value = [double]"2.990417665493794e-06"
$ws.Range("S2").Value = [double]"2.990417665493794e-06"
$ws.Range("T2").Value = 68.32894737885343
$ws.Range("U2").Value = "[51.34761870469515, 85.3102760530117]"
$ws.Range("V2").Value = [double]"2.43157494139723e-10"
$ws.Range("W2").Value = [double]"2.43157494139723e-10"
$ws.Range("X2").Value = 18.46846846846877
$ws.Range("Y2").Value = 16.15375375375402
$ws.Range("Z2").Value = 20.78318318318352
$ws.Range("F3").Value = 24.60000000000041
$ws.Range("H3").Value = 0.000304559793778636
$ws.Range("I3").Value = 0.000304559793778636
$ws.Range("L3").Value = 47.89697221210316
$ws.Range("M3").Value = "[22.360917332098694, 73.43302709210762]"
$ws.Range("N3").Value = 0.000462218280720128
$ws.Range("O3").Value = 0.000462218280720128
$ws.Range("P3").Value = 1.654131867655887
$ws.Range("Q3").Value = "[0.9874475407679624, 2.320816194543812]"
$ws.Range("R3").Value = [double]"9.283516661628965e-06"
$ws.Range("S3").Value = [double]"9.283516661628965e-06"
$ws.Range("T3").Value = 63.23964304542954
$ws.Range("U3").Value = "[47.59079426912699, 78.88849182173209]"
$ws.Range("V3").Value = [double]"2.162963141927321e-10"
$ws.Range("W3").Value = [double]"2.162963141927321e-10"
$ws.Range("X3").Value = 18.12372372372402
$ws.Range("Y3").Value = 15.51351351351377
$ws.Range("Z3").Value = 20.73393393393428
$ws.Range("F4").Value = 24.60000000000041
$ws.Range("H4").Value = 0.0001061576332401781
$ws.Range("I4").Value = 0.0001061576332401781
$ws.Range("L4").Value = 50.4990393350653
$ws.Range("M4").Value = "[23.222127406782008, 77.77595126334859]"
$ws.Range("N4").Value = 0.0005362379002795681
$ws.Range("O4").Value = 0.0005362379002795681
$ws.Range("P4").Value = 1.213868632918579
$ws.Range("Q4").Value = "[0.5849211547224238, 1.8428161111147352]"
$ws.Range("R4").Value = 0.0003306931306712446
$ws.Range("S4").Value = 0.0003306931306712446
$ws.Range("T4").Value = 64.27061226924053
$ws.Range("U4").Value = "[49.03583722023629, 79.50538731824477]"
$ws.Range("V4").Value = [double]"6.59525767332525e-11"
$ws.Range("W4").Value = [double]"6.59525767332525e-11"
$ws.Range("X4").Value = 19.84744744744777
$ws.Range("Y4").Value = 17.38498498498526
$ws.Range("Z4").Value = 22.30990990991028
$ws.Range("F5").Value = 24.60000000000041
$ws.Range("H5").Value = [double]"6.661589142575863e-05"
$ws.Range("I5").Value = [double]"6.661589142575863e-05"
$ws.Range("L5").Value = 63.35104491000297
$ws.Range("M5").Value = "[34.20698428813414, 92.4951055318718]"
$ws.Range("N5").Value = [double]"7.053171709570094e-05"
$ws.Range("O5").Value = [double]"7.053171709570094e-05"
$ws.Range("P5").Value = 1.490605523324888
$ws.Range("Q5").Value = "[0.8868159442565782, 2.0943951023931975]"
$ws.Range("R5").Value = [double]"1.008677190750262e-05"
$ws.Range("S5").Value = [double]"1.008677190750262e-05"
$ws.Range("T5").Value = 68.69178241630826
$ws.Range("U5").Value = "[50.07356656728706, 87.30999826532945]"
$ws.Range("V5").Value = [double]"2.349141103863417e-09"
$ws.Range("W5").Value = [double]"2.349141103863417e-09"
$ws.Range("X5").Value = 18.76396396396427
$ws.Range("Y5").Value = 16.40000000000026
$ws.Range("Z5").Value = 21.12792792792828
$ws.Range("F6").Value = 24.60000000000041
$ws.Range("H6").Value = [double]"1.324943471159301e-05"
$ws.Range("I6").Value = [double]"1.324943471159301e-05"
$ws.Range("L6").Value = 57.21879793416412
$ws.Range("M6").Value = "[30.22380596528069, 84.21378990304756]"
$ws.Range("N6").Value = [double]"9.995088601266744e-05"
$ws.Range("O6").Value = [double]"9.995088601266744e-05"
$ws.Range("P6").Value = 1.075500187715424
$ws.Range("Q6").Value = "[0.5597632555945777, 1.5912371198362711]"
$ws.Range("R6").Value = 0.0001244302187475288
$ws.Range("S6").Value = 0.0001244302187475288
$ws.Range("T6").Value = 59.05883754966263
$ws.Range("U6").Value = "[43.912825803372996, 74.20484929595227]"
$ws.Range("V6").Value = [double]"5.635036881557198e-10"
$ws.Range("W6").Value = [double]"5.635036881557198e-10"
$ws.Range("X6").Value = 20.38918918918953
$ws.Range("Y6").Value = 18.36996996997027
$ws.Range("Z6").Value = 22.40840840840878
$ws.Range("B7").Value = 0
$ws.Range("F7").Value = 24.60000000000041
$ws.Range("H7").Value = 0.01189924204556492
$ws.Range("I7").Value = 0.01189924204556492
$ws.Range("L7").Value = 37.79035396836994
$ws.Range("M7").Value = "[5.413561756656378, 70.1671461800835]"
$ws.Range("N7").Value = 0.02317089288987062
$ws.Range("O7").Value = 0.02317089288987062
$ws.Range("P7").Value = 0.7358685494895019
$ws.Range("Q7").Value = "[-0.018868424345884094, 1.4906055233248878]"
$ws.Range("R7").Value = 0.05575518630720189
$ws.Range("S7").Value = 0.05575518630720189
$ws.Range("T7").Value = 61.40480430099272
$ws.Range("U7").Value = "[44.295349923972196, 78.51425867801325]"
$ws.Range("V7").Value = [double]"4.675551723565263e-09"
$ws.Range("W7").Value = [double]"4.675551723565263e-09"
$ws.Range("X7").Value = 21.71891891891928
$ws.Range("Y7").Value = 18.76396396396427
$ws.Range("Z7").Value = 24.67387387387428
$ws.Range("F8").Value = 23.87000000000029
$ws.Range("H8").Value = [double]"2.821498790750443e-05"
$ws.Range("I8").Value = [double]"2.821498790750443e-05"
$ws.Range("L8").Value = 64.70796341730717
$ws.Range("M8").Value = "[29.99790565497122, 99.41802117964312]"
$ws.Range("N8").Value = 0.0004956632778396752
$ws.Range("O8").Value = 0.0004956632778396752
$ws.Range("P8").Value = 0.748447499053424
$ws.Range("Q8").Value = "[0.24528951649650033, 1.2516054816103477]"
$ws.Range("R8").Value = 0.004437961677769886
$ws.Range("S8").Value = 0.004437961677769886
$ws.Range("T8").Value = 67.64805933484189
$ws.Range("U8").Value = "[49.83030889323723, 85.46580977644655]"
$ws.Range("V8").Value = [double]"1.131118976260836e-09"
$ws.Range("W8").Value = [double]"1.131118976260836e-09"
$ws.Range("X8").Value = 21.02662662662689
$ws.Range("Y8").Value = 19.11511511511535
$ws.Range("Z8").Value = 22.93813813813842
$ws.Range("F9").Value = 23.87000000000029
$ws.Range("H9").Value = 0.00458621399561121
$ws.Range("I9").Value = 0.00458621399561121
$ws.Range("L9").Value = 43.06300640838877
$ws.Range("M9").Value = "[10.523371465597734, 75.60264135117981]"
$ws.Range("N9").Value = 0.01063867467344215
$ws.Range("O9").Value = 0.01063867467344215
$ws.Range("P9").Value = 1.150973885098963
$ws.Range("Q9").Value = "[0.3459211130078854, 1.9560266571900398]"
$ws.Range("R9").Value = 0.006076564507542459
$ws.Range("S9").Value = 0.006076564507542459
$ws.Range("T9").Value = 71.43759914686389
$ws.Range("U9").Value = "[53.89453409246363, 88.98066420126415]"
$ws.Range("V9").Value = [double]"1.756761403015616e-10"
$ws.Range("W9").Value = [double]"1.756761403015616e-10"
$ws.Range("X9").Value = 19.49741741741766
$ws.Range("Y9").Value = 16.43899899899921
$ws.Range("Z9").Value = 22.55583583583611
$ws.Range("F10").Value = 23.87000000000029
$ws.Range("H10").Value = [double]"8.038386067332759e-06"
$ws.Range("I10").Value = [double]"8.038386067332759e-06"
$ws.Range("L10").Value = 53.70505883731979
$ws.Range("M10").Value = "[28.420968153355133, 78.98914952128445]"
$ws.Range("N10").Value = [double]"9.713221821816553e-05"
$ws.Range("O10").Value = [double]"9.713221821816553e-05"
$ws.Range("P10").Value = 1.150973885098963
$ws.Range("Q10").Value = "[0.6478159025420389, 1.6541318676558863]"
$ws.Range("R10").Value = [double]"3.358797948682657e-05"
$ws.Range("S10").Value = [double]"3.358797948682657e-05"
$ws.Range("T10").Value = 57.22369511901586
$ws.Range("U10").Value = "[43.472488629277336, 70.97490160875438]"
$ws.Range("V10").Value = [double]"9.665335198860703e-11"
$ws.Range("W10").Value = [double]"9.665335198860703e-11"
$ws.Range("X10").Value = 19.49741741741766
$ws.Range("Y10").Value = 17.58590590590612
$ws.Range("Z10").Value = 21.40892892892919
$ws.Range("F11").Value = 23.87000000000029
$ws.Range("H11").Value = 0.02239038194885656
$ws.Range("I11").Value = 0.02239038194885656
$ws.Range("L11").Value = 37.70962419135108
$ws.Range("M11").Value = "[4.729006043150065, 70.6902423395521]"
$ws.Range("N11").Value = 0.02596213138206882
$ws.Range("O11").Value = 0.02596213138206882
$ws.Range("P11").Value = 1.264184431174272
$ws.Range("Q11").Value = "[0.081763172165501, 2.446605690183042]"
$ws.Range("R11").Value = 0.03668453296815799
$ws.Range("S11").Value = 0.03668453296815799
$ws.Range("T11").Value = 70.07440124560506
$ws.Range("U11").Value = "[51.448455520980644, 88.70034697022948]"
$ws.Range("V11").Value = [double]"1.430398466339966e-09"
$ws.Range("W11").Value = [double]"1.430398466339966e-09"
$ws.Range("X11").Value = 19.06732732732756
$ws.Range("Y11").Value = 14.57527527527546
$ws.Range("Z11").Value = 23.55937937937967
$ws.Range("F12").Value = 23.87000000000029
$ws.Range("H12").Value = [double]"8.401675508407092e-05"
$ws.Range("I12").Value = [double]"8.401675508407092e-05"
$ws.Range("L12").Value = 57.89314637332472
$ws.Range("M12").Value = "[24.16041221653896, 91.62588053011048]"
$ws.Range("N12").Value = 0.001205383077964806
$ws.Range("O12").Value = 0.001205383077964806
$ws.Range("P12").Value = 0.7610264486173479
$ws.Range("Q12").Value = "[0.24528951649649944, 1.2767633807381964]"
$ws.Range("R12").Value = 0.00473678447755721
$ws.Range("S12").Value = 0.00473678447755721
$ws.Range("T12").Value = 61.5750534536878
$ws.Range("U12").Value = "[44.48838468100131, 78.6617222263743]"
$ws.Range("V12").Value = [double]"4.22594093052453e-09"
$ws.Range("W12").Value = [double]"6.59525767332525e-11"
$ws.Range("X12").Value = 20.97883883883909
$ws.Range("Y12").Value = 19.01953953953976
$ws.Range("Z12").Value = 22.93813813813842
